# Add a new worksheet named "Sheet2" right after "Sheet1" and use it
# as a log of thoughts (ID / Date / Log columns).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert the new sheet right after Sheet1 so it keeps Sheet1 in first
# position and becomes the (now active) second tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "Date"
$ws2.Range("C1").Value = "Log"

# Widen the Log column so entries are readable
$ws2.Columns.Item(3).ColumnWidth = 104.7109375

# Leave the selection ready on the first data row, as in the original edit
$ws2.Range("A2").Select()
